$wb = $excel.ActiveWorkbook

# Sheet 1 currently holds hotel_info data (with header + 1 data row, 9 cols)
# Sheet 2 currently holds review_info headers only (1 row, 25 cols)
$sheetHotel = $wb.Worksheets.Item(1)
$sheetReview = $wb.Worksheets.Item(2)

# Capture the existing review_info headers (25 columns) before clearing anything
$reviewHeaders = @()
for ($c = 1; $c -le 25; $c++) {
    $reviewHeaders += $sheetReview.Cells.Item(1, $c).Value2
}

# Capture the existing hotel_info headers (9 columns) and data row
$hotelHeaders = @()
for ($c = 1; $c -le 9; $c++) {
    $hotelHeaders += $sheetHotel.Cells.Item(1, $c).Value2
}
$hotelData = @()
for ($c = 1; $c -le 9; $c++) {
    $hotelData += $sheetHotel.Cells.Item(2, $c).Value2
}

# Clear both sheets entirely
$sheetHotel.Cells.Clear()
$sheetReview.Cells.Clear()

# Rename via temporary names first to avoid name collisions during the swap
$sheetHotel.Name = "TempSheetA"
$sheetReview.Name = "TempSheetB"

# --- Rewrite sheet that used to be "hotel_info" (position 1) as "review_info" ---
$sheetHotel.Name = "review_info"
for ($c = 1; $c -le 25; $c++) {
    $sheetHotel.Cells.Item(1, $c).Value = $reviewHeaders[$c - 1]
}

# --- Rewrite sheet that used to be "review_info" (position 2) as "hotel_info" ---
# Original columns (A..I): STR, Hotel_Name, City, Zip, TA_ReviewURL,
#   Tripadvisor_Hotel_Name, English_Reviews_num, Local_Rank, Total_Reviews_num
# New columns (A..J): STR, Hotel_Name, State, City, Zip, TA_ReviewURL,
#   Tripadvisor_Hotel_Name, English_Reviews_num, Local_Rank, Total_Reviews_num
$sheetReview.Name = "hotel_info"

$newHotelHeaders = @($hotelHeaders[0], $hotelHeaders[1], "State", $hotelHeaders[2], $hotelHeaders[3], $hotelHeaders[4], $hotelHeaders[5], $hotelHeaders[6], $hotelHeaders[7], $hotelHeaders[8])
$newHotelData = @($hotelData[0], $hotelData[1], "Louisiana", $hotelData[2], $hotelData[3], $hotelData[4], $hotelData[5], $hotelData[6], $hotelData[7], $hotelData[8])

for ($c = 1; $c -le $newHotelHeaders.Length; $c++) {
    $sheetReview.Cells.Item(1, $c).Value = $newHotelHeaders[$c - 1]
}

# Columns H, I, J (English_Reviews_num, Local_Rank, Total_Reviews_num) hold
# numeric-looking text in the source data ("819", "10", "823") and must be
# written back as text, not numbers, so force a text number format first.
$textCols = @(8, 9, 10)
foreach ($tc in $textCols) {
    $sheetReview.Cells.Item(2, $tc).NumberFormat = "@"
}

for ($c = 1; $c -le $newHotelData.Length; $c++) {
    $sheetReview.Cells.Item(2, $c).Value = $newHotelData[$c - 1]
}

# Re-assert exact types/values for the numeric-looking cells, since reading
# them back through the object model can silently coerce types:
#  - STR id and Zip must be true numbers
#  - English_Reviews_num / Local_Rank / Total_Reviews_num must stay text
$sheetReview.Cells.Item(2, 1).Value = 62888
$sheetReview.Cells.Item(2, 5).Value = 70461
$sheetReview.Cells.Item(2, 8).Value = "819"
$sheetReview.Cells.Item(2, 9).Value = "10"
$sheetReview.Cells.Item(2, 10).Value = "823"

# Drop back to the default "Normal" style so the forced text format doesn't
# leave a stray cell style behind (value/type are already locked in above).
foreach ($tc in $textCols) {
    $sheetReview.Cells.Item(2, $tc).Style = "Normal"
}
